# Applies the LOM3205.docx edit:
#   1. Bump the "Ativação" date from 2019 to 2023.
#   2/3/4. Add the missing italicised English translation paragraphs
#      under "Objetivos", "Programa resumido" and "Programa".
#   5. Rewrite/expand the Portuguese "Programa" paragraph text.

$d = $word.ActiveDocument

# --- 1. Update activation date (simple in-place text swap) ---
[void]$d.Content.Find.Execute("Ativação: 01/01/2019", $false, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2023", 2)

# --- 2. English translation after the "Objetivos" body paragraph ---
$rObjetivos = $d.Content
$rObjetivos.Find.Text = "Estudo formal da teoria dos campos eletromagnéticos independentes do tempo ou para situações quase-estáticas. Teoria das ondas eletromagnéticas."
[void]$rObjetivos.Find.Execute()
$pObjetivos = $rObjetivos.Paragraphs(1)
$pObjetivos.Range.InsertParagraphAfter()
$rObjetivosEn = $pObjetivos.Next().Range
[void]$rObjetivosEn.MoveEnd(1, -1)
$rObjetivosEn.Text = "Formal study of the theory of time-independent electromagnetic fields or for near static situations. Electromagnetic waves theory."
$rObjetivosEn.Font.Italic = $true

# --- 3. English translation after the "Programa resumido" body paragraph ---
$rResumido = $d.Content
$rResumido.Find.Text = "Eletrostática. Magnetostática. Campos variantes no tempo. Equações de Maxwell. Ondas eletromagnéticas."
[void]$rResumido.Find.Execute()
$pResumido = $rResumido.Paragraphs(1)
$pResumido.Range.InsertParagraphAfter()
$rResumidoEn = $pResumido.Next().Range
[void]$rResumidoEn.MoveEnd(1, -1)
$rResumidoEn.Text = "Electrostatics. Magnetostatic. Time-varying fields. Maxwell's equations. Electromagnetic waves"
$rResumidoEn.Font.Italic = $true

# --- 4. Rewrite the "Programa" paragraph text, then add its English translation ---
$rPrograma = $d.Content
$rPrograma.Find.Text = "Eletrostática (Campo Eletrostático; Potencial Elétrico; Trabalho e Energia em Eletrostática)  Técnicas Especiais (Equações de Laplace; Método das imagens; Separação de variáveis e Expansão em Multipolos) Campo Elétrico da Matéria (Polarização elétrica; Campo de objeto polarizado; cargas ligadas; deslocamento elétrico; Dielétricos lineares) Magnetostática (Lei de Lorentz; Lei de Biot-Savart; Lei de Ampére; Potencial Vetor Magnético) Campo Magnético na Matéria (Magnetização; Campos de objeto magnetizado; Campo auxiliar H; Meios Lineares e não lineares) Eletrodinâmica (Força eletromotriz; Indução eletromagnética; Equações de Maxwell; Leis de conservação) Ondas eletromagnéticas (Propagação no vácuo e na matéria; Reflexão e transmissão) ou Equação de ondas (Planas, esféricas e cilíndricas) e condições de contorno (interfaces)"
[void]$rPrograma.Find.Execute()
$pPrograma = $rPrograma.Paragraphs(1)
$pPrograma.Range.Text = "Eletrostática (campo eletrostático; potencial elétrico; trabalho e energia em eletrostática).  Técnicas especiais para a resolução da equação de Laplace (método das imagens; separação de variáveis). Campo elétrico da matéria (polarização elétrica; campo de objeto polarizado; cargas ligadas; deslocamento elétrico; dielétricos (lineares). Magnetostática (Lei de Lorentz; Lei de Biot-Savart; Lei de Ampére; vetor potencial magnético).  Campo magnético na matéria (magnetização; campos de objeto magnetizado; campo auxiliar H; Eletrodinâmica (força eletromotriz; indução eletromagnética; equações de Maxwell; lei de conservação de carga). Ondas eletromagnéticas (propagação no vácuo e na matéria; reflexão e transmissão), equação de ondas (planas)  e condições de contorno (interfaces). Radiação de dipolo elétrico."
$pPrograma.Range.InsertParagraphAfter()
$rProgramaEn = $pPrograma.Next().Range
[void]$rProgramaEn.MoveEnd(1, -1)
$rProgramaEn.Text = "Electrostatics (electrostatic field; electric potential; work and energy in electrostatics). Special techniques for solving the Laplace’s equation (method of images; separation of variables). Electric field in matter (electric polarization; polarized object field; bound charges; electric displacement; dielectrics (linear). Magnetostatics (Lorentz's law; Biot-Savart's law; Ampere's law; vector magnetic potential). Magnetic field in matter (magnetization; field of a magnetized object; auxiliary field H). Electrodynamics (electromotive force; electromagnetic induction; Maxwell's equations; law of conservation of charge). Electromagnetic waves (propagation in vacuum and in matter; reflection and transmission), plane wave equation and boundary conditions (interfaces). Electric dipole radiation."
$rProgramaEn.Font.Italic = $true
